$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$vins = @(
    "1HGCM82633A123456",
    "JH4KA8260MC012345",
    "1FTSW21R08ED56789",
    "4T1BG22KXWU123456",
    "2HGFG12856H654321",
    "1J4GL48KX7W543210",
    "1N4AL11D75C123456",
    "2T2BZMCA2GC123456",
    "5N1AT2MV5HC789012",
    "JHMFA16586S012345",
    "1HGCM66586A654321",
    "4S4BP61C667789012",
    "1C4NJPFA8CD123456",
    "5XYKTDA26CG654321",
    "1FMEU15N1SLC12345",
    "3C4PDCBG0FT654321",
    "1G1JC1244VJ123456",
    "2G4WB52K8X1234567",
    "2HGFA16578H654321",
    "1HGEJ8240WL123456",
    "JTHBE1BL5D5001234",
    "1FTEX1C85AF123456",
    "1C6RR7MT4FS654321",
    "1FTFW1ET2EFB12345",
    "3FAHP08148R654321",
    "5J6RE4H56BL123456",
    "1HGCP26798A123456",
    "JTLKE50E591654321",
    "1G4HR52KXWH123456",
    "1FMRU1763VLA12345",
    "JHMEG8650PS654321",
    "JN8AZ1MU1AW123456",
    "1FTNW21PX7EC12345",
    "2LMPJ8LR8KBL12345",
    "1FADP3F24EL123456",
    "4S3BMBC69H1254321",
    "1FTRX18L7YKC12345",
    "1HGFA16527L654321",
    "3FA6P0HR5DR123456",
    "1HGCP36759A123456",
    "JN8DR09Y82W654321",
    "3VWDP7AJ4DM123456",
    "5J6YH28696L123456",
    "5NPEC4AC6BH654321",
    "JNKCA21D5VT123456",
    "2HGFG4A55DH654321",
    "1C3CCCAB3FN123456",
    "1GNDT13S692654321",
    "3CZRU5H36GM123456",
    "JN1CV6EK6HM654321",
    "2HKRM4H56DH123456",
    "3FADP4EJ0FM123456",
    "5YJSA1H26EFP67890",
    "1G1ZD5STXJF654321",
    "1N6BA0CJ4FN123456",
    "WBA3A9C57EF456789",
    "1FTEX1CM4FKE12345",
    "JHMGE8H56DC123456",
    "2HGFC2F59HH654321",
    "3VWLL7AJ7EM123456",
    "5TDZK3EH0FS987654",
    "2G1WB58K489123456",
    "1GKS2BKC6GR123456",
    "JN8AS5MT9CW987654",
    "KMHCT4AE6HU123456",
    "1D7RV1CP1AS654321",
    "3GTU2NEC1EG123456",
    "1FTFW1EF3EK456789",
    "2C4RDGBG3FR123456",
    "5XYPGDA35GG654321",
    "1FAHP3FN1AW654321",
    "5FNYF4H99BB123456",
    "JA32U2FU6AU987654",
    "3N1BC1AP8CL123456",
    "5XXGN4A76EG123456",
    "KM8SRDHF9EU654321",
    "3VW517AT9DM456789",
    "JTHBK262875123456",
    "1HGEJ622XWL654321",
    "2C3CCARG6FH123456",
    "KM8K53A14JU123456",
    "1C4PJMDS5GW123456",
    "JN8AF5MR5ET654321",
    "5TFAW5F12DX987654",
    "4T4BF3EK9BR654321",
    "3C6UR5FL1FG123456",
    "1GCGC24R1ES654321",
    "1GTW7AFG6G1212345",
    "5FRYD3H29GB654321",
    "5TFHW5F17GX123456",
    "JTEBU5JR3E5678901",
    "3FA6P0LU5HR654321",
    "1GCGK29J1LE123456",
    "1GNER13DXYJ654321",
    "JHMCM56457S123456",
    "4T1BE32K33U123456",
    "2G1WT57K391234567",
    "5J8TB4H36HL654321",
    "JN8AZ2KR6HT654321",
    "KMHCM3ACXBU123456",
    "1HGCR2F32EA987654"
)

for ($i = 0; $i -lt $vins.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 13).Value = $vins[$i]
}

$ws.Range("M111").Select()

